# Gallery_Accessories.xlsx - add a new "Slovakia" market sheet
# ---------------------------------------------------------------
# Mirrors what a user would do in Excel: duplicate an existing
# "market" sheet (Germany - it already has the same row heights/
# layout the new sheet needs), rename it, trim two accessory rows
# that don't apply to the Slovakia market, fill in the
# market-specific values, resize column B to fit the new text and
# finally leave the new sheet selected/active - exactly like
# Excel does right after you finish editing a freshly added sheet.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("Germany")
$portugal = $wb.Worksheets.Item("Portugal")

# Duplicate the template sheet and drop it right after "Portugal"
# (i.e. at the end of the tab strip).
$template.Copy($null, $portugal)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Slovakia"

# This market's accessory list skips "Mounting Frame" and
# "Blank CUI Inserts" - remove those two rows.
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(11).Delete()

# Fill in the market-specific code / label. Order matters here: the
# workbook's shared-string table is built in the order new strings
# are encountered, and the target file expects
# "NGC-2930/T3219/T3183" before "Slovakia Market".
$ws.Range("B4").Value = "NGC-2930/T3219/T3183"
$ws.Range("B2").Value = "Slovakia Market"

# Column A/D keep the same widths as the other market sheets; column
# B needs to widen so the new, longer market name/code fit ("best
# fit").
$ws.Columns.Item(1).ColumnWidth = 25.5546875
$ws.Columns.Item(4).ColumnWidth = 17.6640625
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(2).ColumnWidth = 20.3

# Leave the new sheet as the active one with its last-used cell
# selected, same as Excel would after you stopped editing it.
$ws.Range("C15").Select()
$ws.Activate()

# "Portugal" is no longer the active tab, so its lingering selection
# becomes the whole-sheet selection left over from Select All.
$portugal.Range("A1:XFD1048576").Select()

$ws.Activate()
